{"js": "// Update the CLS NOFO eligibility paragraph:\n//   - keep the leading \". \"\n//   - replace the old eligibility text with the new wording.\nconst oldText =\n  \"In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services. Eligible applicants may include public agencies and nonprofit agencies that provide other services, but eligible applicants must include legal services in their core services. Eligible applicants must also demonstrate a record of providing effective direct services to crime victims.\";\n\nconst newText =\n  \"In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services through the use of attorneys. If a victim service agency doesn\\u2019t currently focus on the provision of legal services, then it is not eligible. \";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\nresults.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Update the CLS NOFO eligibility paragraph wording.\n# The old run: \". In brief, ... crime victims.\" is split into two runs with\n# identical run formatting: \". \" and the new eligibility sentence(s).\n$d = $word.ActiveDocument\n\n$oldText = \". In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services. Eligible applicants may include public agencies and nonprofit agencies that provide other services, but eligible applicants must include legal services in their core services. Eligible applicants must also demonstrate a record of providing effective direct services to crime victims.\"\n\n# Locate the exact run via Find, scoped so the match never extends into the\n# preceding \"NOFO\" run (avoids the engine's run-coalescing on plain edits).\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Target sentence not found in document.\"\n}\n\n# Re-materialize a plain Range over the matched span; InsertXML needs a\n# freshly-built Range (not the live Find range) to replace in place instead\n# of leaving the old text and appending after it.\n$target = $d.Range($searchRange.Start, $searchRange.End)\n\n$run2Text = \"In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services through the use of attorneys. If a victim service agency doesn\" + [char]0x2019 + \"t currently focus on the provision of legal services, then it is not eligible. \"\n\n# Replace the matched range's content with explicit OOXML so the two new\n# runs keep the original run formatting, and the preceding \"NOFO\" run is\n# left completely untouched (InsertXML only rewrites the addressed range).\n$xml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:color w:val=\"201F1E\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">. </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:color w:val=\"201F1E\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">$run2Text</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$target.InsertXML($xml)\n"}
